$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1: add two new header cells (V1, W1), matching the header style
# (bold white text on dark fill) used by the rest of row 1 ---
$ws.Range("V1").Value = "Client ID:*"
$ws.Range("W1").Value = "FORMAT ID:*"
$ws.Range("A1").Copy()
$ws.Range("V1:W1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2: remove the old sample/demo data (A2:T2) so the row becomes
# blank like rows 3-5, keeping only a styled placeholder in B2 ---
$ws.Range("A2:T2").Clear()
$ws.Rows.Item(2).AutoFit()
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remove the trailing rows 6 and 7 so the sheet only spans to row 5 ---
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# --- Restore the selection to match the saved workbook state ---
$ws.Range("B18").Select()

$wb.Save()
